$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Suite")
$ws.Range("B4").Value = "N"
